# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# Sheet "展览" (sheet1): row 3 -> event "合肥·第九届环形宇宙动漫游戏嘉年华" (F3: 3230 -> 3243)
#                        row 6 -> event "合肥·心动恋章·冬日序国乙&代号鸢同人only" (F6: 139 -> 144)
# Sheet "全部类型" (sheet4, aggregated list): same two events appear again at
#                        row 7 (F7: 3230 -> 3243) and row 11 (F11: 139 -> 144)

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 3243
$wsExpo.Range("F6").Value = 144

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 3243
$wsAll.Range("F11").Value = 144
